$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for the PR #1420 tests (sc17, sc18, sc19)
$ws.Range("A21").Value = "sc17"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = "Hierarchal conditions. All locations duplicated some with blank cond tag. Tests PR#1420"
$ws.Range("F21").Value = "complete"
$ws.Range("G21").Value = "yes"
$ws.Range("H21").Value = "done"

$ws.Range("A22").Value = "sc18"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = "Simple two location version of sc17. Includes duplicate loc 2 with blank cond tag."
$ws.Range("F22").Value = "input files"
$ws.Range("G22").Value = "no"
$ws.Range("H22").Value = "to do"

$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = "Simple two location version of sc17. No duplicate loc 2."
$ws.Range("A23").Value = "sc19"
$ws.Range("F23").Value = "input files"
$ws.Range("G23").Value = "no"
$ws.Range("H23").Value = "to do"

# Update the selection to match the committed workbook state
$ws.Range("A22:H23").Select()
